$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 (I0) and J1 (IF), copying the bold/border header style from H1
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Populate I2:J82 with the I0/IF numeric data for each row
$arr = New-Object "object[,]" 81,2
$arr[0,0] = 9
$arr[0,1] = 9
$arr[1,0] = 3
$arr[1,1] = 3
$arr[2,0] = 5
$arr[2,1] = 6
$arr[3,0] = 3
$arr[3,1] = 3
$arr[4,0] = 2
$arr[4,1] = 3
$arr[5,0] = 7
$arr[5,1] = 8
$arr[6,0] = 3
$arr[6,1] = 3
$arr[7,0] = 5
$arr[7,1] = 5
$arr[8,0] = 3
$arr[8,1] = 4
$arr[9,0] = 5
$arr[9,1] = 5
$arr[10,0] = 8
$arr[10,1] = 9
$arr[11,0] = 6
$arr[11,1] = 6
$arr[12,0] = 7
$arr[12,1] = 7
$arr[13,0] = 10
$arr[13,1] = 11
$arr[14,0] = 10
$arr[14,1] = 10
$arr[15,0] = 6
$arr[15,1] = 7
$arr[16,0] = 6
$arr[16,1] = 6
$arr[17,0] = 9
$arr[17,1] = 9
$arr[18,0] = 7
$arr[18,1] = 7
$arr[19,0] = 6
$arr[19,1] = 6
$arr[20,0] = 5
$arr[20,1] = 5
$arr[21,0] = 7
$arr[21,1] = 7
$arr[22,0] = 6
$arr[22,1] = 6
$arr[23,0] = 7
$arr[23,1] = 7
$arr[24,0] = 7
$arr[24,1] = 7
$arr[25,0] = 9
$arr[25,1] = 9
$arr[26,0] = 8
$arr[26,1] = 8
$arr[27,0] = 9
$arr[27,1] = 9
$arr[28,0] = 8
$arr[28,1] = 8
$arr[29,0] = 8
$arr[29,1] = 9
$arr[30,0] = 7
$arr[30,1] = 7
$arr[31,0] = 4
$arr[31,1] = 4
$arr[32,0] = 8
$arr[32,1] = 8
$arr[33,0] = 9
$arr[33,1] = 9
$arr[34,0] = 10
$arr[34,1] = 10
$arr[35,0] = 10
$arr[35,1] = 10
$arr[36,0] = 9
$arr[36,1] = 9
$arr[37,0] = 9
$arr[37,1] = 9
$arr[38,0] = 9
$arr[38,1] = 9
$arr[39,0] = 9
$arr[39,1] = 9
$arr[40,0] = 9
$arr[40,1] = 9
$arr[41,0] = 9
$arr[41,1] = 9
$arr[42,0] = 9
$arr[42,1] = 9
$arr[43,0] = 8
$arr[43,1] = 8
$arr[44,0] = 9
$arr[44,1] = 9
$arr[45,0] = 9
$arr[45,1] = 9
$arr[46,0] = 9
$arr[46,1] = 9
$arr[47,0] = 9
$arr[47,1] = 9
$arr[48,0] = 10
$arr[48,1] = 10
$arr[49,0] = 9
$arr[49,1] = 9
$arr[50,0] = 9
$arr[50,1] = 9
$arr[51,0] = 9
$arr[51,1] = 9
$arr[52,0] = 9
$arr[52,1] = 9
$arr[53,0] = 9
$arr[53,1] = 9
$arr[54,0] = 8
$arr[54,1] = 9
$arr[55,0] = 9
$arr[55,1] = 9
$arr[56,0] = 9
$arr[56,1] = 9
$arr[57,0] = 9
$arr[57,1] = 9
$arr[58,0] = 9
$arr[58,1] = 9
$arr[59,0] = 10
$arr[59,1] = 10
$arr[60,0] = 9
$arr[60,1] = 9
$arr[61,0] = 9
$arr[61,1] = 9
$arr[62,0] = 9
$arr[62,1] = 9
$arr[63,0] = 10
$arr[63,1] = 10
$arr[64,0] = 9
$arr[64,1] = 9
$arr[65,0] = 10
$arr[65,1] = 10
$arr[66,0] = 9
$arr[66,1] = 10
$arr[67,0] = 9
$arr[67,1] = 9
$arr[68,0] = 8
$arr[68,1] = 8
$arr[69,0] = 8
$arr[69,1] = 8
$arr[70,0] = 8
$arr[70,1] = 8
$arr[71,0] = 8
$arr[71,1] = 8
$arr[72,0] = 9
$arr[72,1] = 9
$arr[73,0] = 8
$arr[73,1] = 8
$arr[74,0] = 7
$arr[74,1] = 7
$arr[75,0] = 8
$arr[75,1] = 8
$arr[76,0] = 8
$arr[76,1] = 8
$arr[77,0] = 5
$arr[77,1] = 5
$arr[78,0] = 8
$arr[78,1] = 8
$arr[79,0] = 4
$arr[79,1] = 4
$arr[80,0] = 4
$arr[80,1] = 4

$ws.Range("I2:J82").Value = $arr
